$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A6").Value = 29208154
$ws.Range("B6").Value = "iamnumber10"
$ws.Range("C6").Value = "hamode"
$ws.Range("D6").Value = "badarni"
$ws.Range("E6").Value = "hamode@gmail.com"
$ws.Range("F6").Value = 548177167
